# [#134279773] Importacao de usuarios aceita informacao de turma
#
# Adds a new "Turma" (class) column (J) to the modelo-importacao sheet,
# filling in "Turma A" / "Turma B" for the first two sample rows (the
# third sample row is left without a turma, matching the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Turma" header + sample data in column J.
$ws.Range("J1").Value = "Turma"
$ws.Range("J2").Value = "Turma A"
$ws.Range("J3").Value = "Turma B"
# Row 4 (Luana Barbosa Dias) intentionally has no turma assigned.

# The existing columns were very slightly narrowed (a cosmetic re-layout
# side effect of adding the new column in the original workbook).
$ws.Columns.Item(1).ColumnWidth = 19.008503401360567
$ws.Columns.Item(2).ColumnWidth = 12.799319727891167
$ws.Columns.Item(3).ColumnWidth = 6.319727891156467
$ws.Columns.Item(4).ColumnWidth = 23.73299319727887
$ws.Columns.Item(5).ColumnWidth = 30.350340136054466
$ws.Columns.Item(7).ColumnWidth = 6.187074829931977
$ws.Columns.Item(8).ColumnWidth = 7.671768707482998
$ws.Columns.Item(9).ColumnWidth = 8.753401360544217

# Tab ratio (horizontal scrollbar / sheet tabs split) nudged slightly.
$excel.ActiveWindow.TabRatio = 0.985

# Selection moved to the newly added cell.
$ws.Range("J4").Select()
